$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 20 de Julio de 2020 a las 00:40"

# Row 4
$ws.Range("B4").Value = 3889130
$ws.Range("C4").Value = 55859
$ws.Range("D4").Value = 1798932
$ws.Range("E4").Value = 1946962
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 359
$ws.Range("H4").Value = 143236

# Row 5
$ws.Range("B5").Value = 2098389
$ws.Range("C5").Value = 23143
$ws.Range("D5").Value = 1371229
$ws.Range("E5").Value = 647672
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 671
$ws.Range("H5").Value = 79488

# Row 9
$ws.Range("B9").Value = 353590
$ws.Range("C9").Value = 4090
$ws.Range("D9").Value = 241955
$ws.Range("E9").Value = 98448
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 189
$ws.Range("H9").Value = 13187

# Row 21
$ws.Range("B21").Value = 197278
$ws.Range("C21").Value = 6578
$ws.Range("D21").Value = 91793
$ws.Range("E21").Value = 98749
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 220
$ws.Range("H21").Value = 6736

# Row 23
$ws.Range("B23").Value = 126755
$ws.Range("C23").Value = 4231
$ws.Range("D23").Value = 54105
$ws.Range("E23").Value = 70390
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 40
$ws.Range("H23").Value = 2260

# Row 24
$ws.Range("B24").Value = 110338
$ws.Range("C24").Value = 339
$ws.Range("D24").Value = 97026
$ws.Range("E24").Value = 4460
$ws.Range("F24").Value = 0

# Row 50: A50 -> Nigeria
$ws.Range("A50").Value = "Nigeria"
$ws.Range("B50").Value = 36663
$ws.Range("C50").Value = 556
$ws.Range("D50").Value = 15105
$ws.Range("E50").Value = 20769
$ws.Range("F50").Value = 0
$ws.Range("G50").Value = 11
$ws.Range("H50").Value = 789

# Row 51: A51 -> Barein
$ws.Range("A51").Value = "Barein"
$ws.Range("B51").Value = 36422
$ws.Range("C51").Value = 418
$ws.Range("D51").Value = 32372
$ws.Range("E51").Value = 3924
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 2
$ws.Range("H51").Value = 126

# Row 54
$ws.Range("D54").Value = 30300
$ws.Range("E54").Value = 1322
$ws.Range("F54").Value = 0

# Row 70
$ws.Range("B70").Value = 13945
$ws.Range("C70").Value = 90
$ws.Range("D70").Value = 8761
$ws.Range("E70").Value = 4825
$ws.Range("F70").Value = 0

# Row 84
$ws.Range("B84").Value = 8733
$ws.Range("C84").Value = 95
$ws.Range("D84").Value = 4106
$ws.Range("E84").Value = 4327
$ws.Range("F84").Value = 0
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 300

# Row 103
$ws.Range("B103").Value = 3721
$ws.Range("C103").Value = 92
$ws.Range("D103").Value = 1918
$ws.Range("E103").Value = 1772
$ws.Range("F103").Value = 0
$ws.Range("G103").Value = 2
$ws.Range("H103").Value = 31

# Row 130: A130 -> Zimbabue
$ws.Range("A130").Value = "Zimbabue"
$ws.Range("B130").Value = 1611
$ws.Range("C130").Value = 133
$ws.Range("D130").Value = 472
$ws.Range("E130").Value = 1114
$ws.Range("F130").Value = 0
$ws.Range("G130").Value = 0
$ws.Range("H130").Value = 25

# Row 131: A131 -> Yemen
$ws.Range("A131").Value = "Yemen"
$ws.Range("B131").Value = 1606
$ws.Range("C131").Value = 25
$ws.Range("D131").Value = 712
$ws.Range("E131").Value = 449
$ws.Range("F131").Value = 0
$ws.Range("G131").Value = 2
$ws.Range("H131").Value = 445

# Row 132: A132 -> Benin
$ws.Range("A132").Value = "Benin"
$ws.Range("B132").Value = 1602
$ws.Range("C132").Value = 0
$ws.Range("D132").Value = 782
$ws.Range("E132").Value = 789
$ws.Range("F132").Value = 0
$ws.Range("H132").Value = 31

# Row 133: A133 -> Ruanda
$ws.Range("A133").Value = "Ruanda"
$ws.Range("B133").Value = 1582
$ws.Range("C133").Value = 43
$ws.Range("D133").Value = 834
$ws.Range("E133").Value = 743
$ws.Range("F133").Value = 0
$ws.Range("H133").Value = 5

# Row 134: A134 -> Nueva Zelanda
$ws.Range("A134").Value = "Nueva Zelanda"
$ws.Range("B134").Value = 1553
$ws.Range("C134").Value = 3
$ws.Range("D134").Value = 1506
$ws.Range("E134").Value = 25
$ws.Range("F134").Value = 0
$ws.Range("H134").Value = 22

# Row 135: A135 -> Mozambique
$ws.Range("A135").Value = "Mozambique"
$ws.Range("B135").Value = 1491
$ws.Range("C135").Value = 56
$ws.Range("D135").Value = 472
$ws.Range("E135").Value = 1009
$ws.Range("F135").Value = 0
$ws.Range("H135").Value = 10

# Row 178: A178 -> Bahamas
$ws.Range("A178").Value = "Bahamas"
$ws.Range("C178").Value = 15
$ws.Range("D178").Value = 91
$ws.Range("E178").Value = 51
$ws.Range("F178").Value = 0
$ws.Range("H178").Value = 11

# Row 179: A179 -> Bermudas
$ws.Range("A179").Value = "Bermudas"
$ws.Range("B179").Value = 153
$ws.Range("C179").Value = 0
$ws.Range("D179").Value = 137
$ws.Range("E179").Value = 7
$ws.Range("F179").Value = 0
$ws.Range("H179").Value = 9

# Row 180: A180 -> Brunei
$ws.Range("A180").Value = "Brunei"
$ws.Range("B180").Value = 141
$ws.Range("C180").Value = 0
$ws.Range("D180").Value = 138
$ws.Range("E180").Value = 0
$ws.Range("F180").Value = 0
$ws.Range("H180").Value = 3
